$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 values (Response)
$ws.Range("B2").Value = 0.97331189884895597
$ws.Range("C2").Value = 9.7855404994651298
$ws.Range("D2").Value = 4.9203750281142797

# Row 3 values (VM Processing)
$ws.Range("B3").Value = 0.027555473982924999
$ws.Range("C3").Value = 0.0276004301943347
$ws.Range("D3").Value = 0.026525640041961299

# Row 4 values (Network Transmission)
$ws.Range("B4").Value = 0.208883632233589
$ws.Range("C4").Value = 2.4118845871233399
$ws.Range("D4").Value = 1.1969367069761001

# Number formats
$ws.Range("B2:D2").NumberFormat = "0.000"
$ws.Range("B4:D4").NumberFormat = "0.000"
$ws.Range("B3:D3").NumberFormat = "0.00"

# Bad style on columns E:F (all 4 rows)
$ws.Range("E1:F4").Style = "Bad"

# Column B width (target stored width 8.42578125 chars; 7.59 is the closest
# COM ColumnWidth input that this engine's pixel-quantized column model maps
# to the nearest representable stored width of 8.5)
$ws.Columns("B").ColumnWidth = 7.59

# Selection
$ws.Range("G9").Select()
